$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.184.32"
$ws.Range("D3").Value = "2.307.76"
$ws.Range("E3").Value = "  +0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "301.25"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.44"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.81%  "
$ws.Range("E7").Value = "  +0.79%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  +1.73%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.68"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +8.44%  "
$ws.Range("E12").Value = "  +0.77%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "17.77"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +3.70%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.93"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +2.06%  "
$ws.Range("D15").Value = "2.666.04"
$ws.Range("E15").Value = "  +0.21%  "
$ws.Range("D16").Value = "2.321.04"
$ws.Range("E16").Value = "  +1.13%  "
$ws.Range("E17").Value = "  -1.63%  "
$ws.Range("D18").Value = "43.085.12"
$ws.Range("E18").Value = "  +0.18%  "
$ws.Range("E19").Value = "  +9.23%  "
$ws.Range("D20").Value = "0.0₃0907"
$ws.Range("E20").Value = "  +0.47%  "
$ws.Range("E21").Value = "  +1.55%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.00"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.73%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.36"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("E24").Value = "  +7.34%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.12%  "
$ws.Range("E26").Value = "  -0.42%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.12"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +3.02%  "
$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "168.59"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.49%  "
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "34.67"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +2.55%  "
$ws.Range("E30").Value = "  -1.06%  "
$ws.Range("E31").Value = "  +0.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.999"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.09%  "
$ws.Range("E33").Value = "  +1.65%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.67"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +4.57%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.59"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.47%  "
$ws.Range("E36").Value = "  -0.77%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0692"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.50%  "
$ws.Range("E38").Value = "  +1.15%  "
$ws.Range("E39").Value = "  +0.28%  "
$ws.Range("E40").Value = "  -0.68%  "
$ws.Range("E41").Value = "  +0.56%  "
$ws.Range("E42").Value = "  +3.00%  "
$ws.Range("D43").Value = "1.981.57"
$ws.Range("E43").Value = "  -0.73%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.26"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -4.51%  "
$ws.Range("E45").Value = "  +4.36%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.83"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.90%  "
$ws.Range("E47").Value = "  +1.74%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "55.34"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +3.60%  "
$ws.Range("E49").Value = "  +4.09%  "
$ws.Range("D50").Value = "2.531.87"
$ws.Range("E50").Value = "  +0.34%  "
$ws.Range("E51").Value = "  +1.20%  "
